$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("A1").Value = "No"
$ws.Range("B1").Value = "Kode Transaksi"
$ws.Range("C1").Value = "Barang"
$ws.Range("D1").Value = "Harga"
$ws.Range("E1").Value = "Jumlah"

# Remove bold/centered formatting from header row (back to default/general)
$headerRange = $ws.Range("A1:E1")
$headerRange.ClearFormats()

# Detail rows
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "TRX001"
$ws.Range("C2").Value = "Chitato"
$ws.Range("D2").Value = 10000
$ws.Range("E2").Value = 5

$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "TRX001"
$ws.Range("C3").Value = "Beef Slice"
$ws.Range("D3").Value = 25000
$ws.Range("E3").Value = 7

$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "TRX001"
$ws.Range("C4").Value = "Indomilk"
$ws.Range("D4").Value = 6500
$ws.Range("E4").Value = 1

# Update selection to match target (active cell D8)
$ws.Range("D8").Select() | Out-Null
